$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 287, shifting existing rows 287..308 down to 288..309
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical to the surrounding rows for this
# market/product, so copy them down from the (now shifted) row 288.
$ws.Range("A287").Value2 = $ws.Range("A288").Value2
$ws.Range("B287").Value2 = $ws.Range("B288").Value2
$ws.Range("C287").Value2 = $ws.Range("C288").Value2
$ws.Range("E287").Value2 = $ws.Range("E288").Value2
$ws.Range("F287").Value2 = $ws.Range("F288").Value2
$ws.Range("G287").Value2 = $ws.Range("G288").Value2
$ws.Range("H287").Value2 = $ws.Range("H288").Value2
$ws.Range("I287").Value2 = $ws.Range("I288").Value2
$ws.Range("R287").Value2 = $ws.Range("R288").Value2

# New row's own data
$ws.Range("D287").Value2 = 45013
$ws.Range("J287").Value2 = 80
$ws.Range("K287").Value2 = 15000
$ws.Range("L287").Value2 = 15000
$ws.Range("M287").Value2 = 15000
$ws.Range("N287").Value2 = "`$/caja 80 unidades"
$ws.Range("O287").Value2 = "Región del Maule"
$ws.Range("P287").Value2 = 188
$ws.Range("Q287").Value2 = 80
